$wb = $excel.ActiveWorkbook
$wsControl = $wb.Worksheets.Item("Control")
$wsData = $wb.Worksheets.Item("Sheet1")

# --- Sheet2 ("Sheet1") data edits ---

# Fill in previously-empty StdDev (column C) values for the existing rows
$wsData.Range("C3").Value = 1.998
$wsData.Range("C4").Value = 2.6760000000000002
$wsData.Range("C5").Value = 4.7089999999999996
$wsData.Range("C6").Value = 10.265000000000001
$wsData.Range("C7").Value = 19.978999999999999
$wsData.Range("C8").Value = 29.131
$wsData.Range("C9").Value = 36.951000000000001
$wsData.Range("C10").Value = 73.593999999999994
$wsData.Range("C11").Value = 110.577

# Add two new data rows (5000 and 7500 password counts)
$wsData.Range("A12").Value = 5000
$wsData.Range("B12").Formula = '=$B$2*A12'
$wsData.Range("C12").Value = 193.61199999999999

$wsData.Range("A13").Value = 7500
$wsData.Range("B13").Formula = '=$B$2*A13'
$wsData.Range("C13").Value = 294.10899999999998

# Trailing formatted-but-empty cell on row 14, matching column B's number format
$wsData.Range("B14").NumberFormat = "0.000"

# --- View state: selection / frozen-pane scroll position ---
# Sheet1 ("Sheet1" / Data sheet) selection moves to C18 (pane stays at A2)
$wsData.Activate()
$wsData.Range("C18").Select()

# Control sheet selection moves to G3, scrolled back up to A2, and ends as the
# active (tab-selected) sheet
$wsControl.Activate()
$wsControl.Range("G3").Select()
